$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.175.07'
$ws.Range("E2").Value = '  -0.64%  '

$ws.Range("D3").Value = '3.874.42'
$ws.Range("E3").Value = '  -1.00%  '

$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '599.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.42%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.77'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.38%  '

$ws.Range("D7").Value = '3.874.41'
$ws.Range("E7").Value = '  -0.82%  '

$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("E9").Value = '  -0.67%  '

$ws.Range("E10").Value = '  -0.52%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.40'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.78%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.458'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.50%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000248'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.44%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.00'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.30%  '

$ws.Range("D15").Value = '4.523.70'
$ws.Range("E15").Value = '  -1.07%  '

$ws.Range("D16").Value = '3.874.05'
$ws.Range("E16").Value = '  -0.66%  '

$ws.Range("D17").Value = '68.138.25'
$ws.Range("E17").Value = '  -0.90%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.27'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +6.74%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.40'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.35%  '

$ws.Range("E20").Value = '  -1.17%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.85'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.85%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '466.55'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.88%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.731'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.42%  '

$ws.Range("E24").Value = '  -4.39%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.44'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.53%  '

$ws.Range("E26").Value = '  +0.86%  '

$ws.Range("E27").Value = '  +0.82%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.04'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.60%  '

$ws.Range("E29").Value = '  +0.14%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.96'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.92%  '

$ws.Range("D31").Value = '4.021.67'
$ws.Range("E31").Value = '  -1.21%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.75'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.78%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.31'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.20%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '31.26'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.52%  '

$ws.Range("E35").Value = '  +2.16%  '

$ws.Range("D36").Value = '3.846.69'
$ws.Range("E36").Value = '  -0.49%  '

$ws.Range("E37").Value = '  -2.26%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.42'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.61%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.03'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.47%  '

$ws.Range("E40").Value = '  +0.61%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.91'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.28%  '

$ws.Range("E42").Value = '  -0.21%  '

$ws.Range("E43").Value = '  -1.42%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '430.32'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.33%  '

$ws.Range("E45").Value = '  -0.09%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '47.36'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.30%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.54'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.66%  '

$ws.Range("B49").Value = 'FLOKI'
$ws.Range("C49").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000275'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.87%  '

$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '144.10'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.89%  '

$ws.Range("E51").Value = '  +3.45%  '
